$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08498203287993761
$ws.Range("C2").Value = 0.9991939073336056
$ws.Range("D2").Value = 0.2184149706477075
$ws.Range("G2").Value = 0.1177988315665668
$ws.Range("H2").Value = 0.99
